# Update countries & provincias Spain
# - Refresh the COVID data snapshot (new case counts for several countries)
# - A handful of neighbouring countries swap rank (re-sorted by total cases),
#   which moves their whole data row (name + stats) down/up by one position
# - Bump the "last updated" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Straightforward data refreshes (country keeps its row) -----------------

# India (row 6)
$ws.Cells.Item(6, 2).Value = 2330327
$ws.Cells.Item(6, 3).Value = 1922
$ws.Cells.Item(6, 4).Value = 1640021
$ws.Cells.Item(6, 5).Value = 644109
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 9
$ws.Cells.Item(6, 8).Value = 46197

# Pakistan (row 17)
$ws.Cells.Item(17, 2).Value = 285921
$ws.Cells.Item(17, 3).Value = 730
$ws.Cells.Item(17, 4).Value = 263193
$ws.Cells.Item(17, 5).Value = 16599
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 17
$ws.Cells.Item(17, 8).Value = 6129

# Israel (row 33)
$ws.Cells.Item(33, 2).Value = 86959
$ws.Cells.Item(33, 3).Value = 366
$ws.Cells.Item(33, 4).Value = 61576
$ws.Cells.Item(33, 5).Value = 24761
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 622

# Uzbekistan (row 62)
$ws.Cells.Item(62, 2).Value = 32215
$ws.Cells.Item(62, 3).Value = 468
$ws.Cells.Item(62, 4).Value = 24090
$ws.Cells.Item(62, 5).Value = 7917
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 4
$ws.Cells.Item(62, 8).Value = 208

# El Salvador (row 73)
$ws.Cells.Item(73, 2).Value = 21269
$ws.Cells.Item(73, 3).Value = 0
$ws.Cells.Item(73, 4).Value = 9897
$ws.Cells.Item(73, 5).Value = 10795
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 7
$ws.Cells.Item(73, 8).Value = 577

# Tailandia (row 116)
$ws.Cells.Item(116, 2).Value = 3356
$ws.Cells.Item(116, 3).Value = 5
$ws.Cells.Item(116, 4).Value = 3169
$ws.Cells.Item(116, 5).Value = 129
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 58

# Islas Turcas y Caicos (row 180)
$ws.Cells.Item(180, 2).Value = 224
$ws.Cells.Item(180, 3).Value = 8
$ws.Cells.Item(180, 4).Value = 39
$ws.Cells.Item(180, 5).Value = 183
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 2

# --- Rows whose rank swapped with their neighbour (name + stats move) ------

# Armenia / Kirguistan swap (rows 55-56)
$ws.Cells.Item(55, 1).Value = "Kirguistan"
$ws.Cells.Item(55, 2).Value = 40759
$ws.Cells.Item(55, 3).Value = 304
$ws.Cells.Item(55, 4).Value = 32997
$ws.Cells.Item(55, 5).Value = 6278
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 6
$ws.Cells.Item(55, 8).Value = 1484
$ws.Cells.Item(56, 1).Value = "Armenia"
$ws.Cells.Item(56, 2).Value = 40593
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(56, 4).Value = 33157
$ws.Cells.Item(56, 5).Value = 6633
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 803

# Papua Nueva Guinea / San Martin (Parte Holandesa) swap (rows 181-182)
$ws.Cells.Item(181, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(181, 2).Value = 219
$ws.Cells.Item(181, 3).Value = 14
$ws.Cells.Item(181, 4).Value = 102
$ws.Cells.Item(181, 5).Value = 100
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 17
$ws.Cells.Item(182, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(182, 2).Value = 214
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 71
$ws.Cells.Item(182, 5).Value = 140
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 3

# Montserrat / Islas Malvinas swap (rows 213-214)
$ws.Cells.Item(213, 1).Value = "Islas Malvinas"
$ws.Cells.Item(213, 2).Value = 13
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 13
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0
$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 2).Value = 13
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1

# --- Timestamp -------------------------------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 08:12"
